$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 209; this shifts rows 209..258 down to 210..259
# and keeps formatting (style) consistent with the row above.
$ws.Rows.Item(209).Insert()

# Populate the newly inserted row 209 with a new data record (same market/region
# constants as its neighbours, but new date / volume / price figures).
$ws.Cells.Item(209, 1).Value = 4
$ws.Cells.Item(209, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(209, 3).Value = "Los Lagos"
$ws.Cells.Item(209, 4).Value = 44754
$ws.Cells.Item(209, 5).Value = 10
$ws.Cells.Item(209, 6).Value = 100112032
$ws.Cells.Item(209, 7).Value = "Zapallo italiano"
$ws.Cells.Item(209, 8).Value = "Sin especificar"
$ws.Cells.Item(209, 9).Value = "Primera"
$ws.Cells.Item(209, 10).Value = 200
$ws.Cells.Item(209, 11).Value = 15000
$ws.Cells.Item(209, 12).Value = 16000
$ws.Cells.Item(209, 13).Value = 15500
$ws.Cells.Item(209, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(209, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(209, 16).Value = 310
$ws.Cells.Item(209, 17).Value = 50
$ws.Cells.Item(209, 18).Value = "Hortaliza"
